$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1. Append a "]" right after the sentence about hand-editing XML files
#    ("This involves hand editing XML files.  Perhaps a dedicated
#    editing tool would help?"), as its own run.
#    Splitting a paragraph, typing into the new one, then re-joining
#    the paragraph mark keeps the two pieces of text as separate runs
#    (instead of being silently coalesced into one run).
# ---------------------------------------------------------------------
$target = "This involves hand editing XML files.  Perhaps a dedicated editing tool would help?"
$rng = $d.Content
$found = $rng.Find.Execute($target, $true, $false, $false, $false, $false,
                            $true, 1, $false, "", 0)
if ($found) {
    $rng.Collapse(0)
    $rng.InsertParagraphAfter()

    foreach ($p in $d.Paragraphs) {
        if ($p.Range.Text -eq ($target + "`r")) {
            $newPara = $p.Next()
            $newPara.Range.InsertAfter("]")
            break
        }
    }

    foreach ($p in $d.Paragraphs) {
        if ($p.Range.Text -eq ($target + "`r")) {
            $markRange = $d.Range($p.Range.End - 1, $p.Range.End)
            $markRange.Delete()
            break
        }
    }
}

# ---------------------------------------------------------------------
# 2. Add a new bullet "Copying a project to a different server." right
#    after the "Restoring a whole project." bullet (same list style).
# ---------------------------------------------------------------------
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -eq "Restoring a whole project.`r") {
        $p.Range.InsertParagraphAfter()
        break
    }
}
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -eq "Restoring a whole project.`r") {
        $newPara = $p.Next()
        $newPara.Range.InsertAfter("Copying a project to a different server.")
        break
    }
}

# ---------------------------------------------------------------------
# 3. Nudge the first table's first two column widths by one twip each
#    (3959 -> 3958 and 541 -> 542), applied uniformly to every row.
# ---------------------------------------------------------------------
$t = $d.Tables(1)
$t.Columns(1).Width = 197.9   # 3958 twips
$t.Columns(2).Width = 27.1    # 542 twips

# ---------------------------------------------------------------------
# 4. Merge the "<> " / "means that the timestamps ..." runs into a
#    single sentence.
# ---------------------------------------------------------------------
$rng2 = $d.Content
$found2 = $rng2.Find.Execute(
    "<> means that the timestamps are the same, but the test case content is different.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "<> means that the timestamps are the same, but the test case content is different.",
    2)
Write-Host "Done"
